$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Stash pristine formats we will need again after Hyperlinks.Add() mutates
#    the cell style of any range it touches. N10 / C10 are inside the
#    existing A1:N42 dimension, already blank/unused, so borrowing them as
#    scratch cells does not alter <dimension>, <cols> spans, etc. We clear
#    their contents/style back to original at the very end.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()                      # B2 currently carries the "Hyperlink" style (s=13)
$ws.Range("N10").PasteSpecial(-4122)        # xlPasteFormats -> stash style 13 on N10

$ws.Range("C10").Copy()                     # C10 carries the plain data style (s=11)
$ws.Range("Z1").Value = $null                # no-op placeholder (kept for clarity)

# ---------------------------------------------------------------------------
# 2) Plain value edits (cells that are not part of a hyperlink target)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "RegistrationTest"
$ws.Range("E2").Value = "Ruby"
$ws.Range("E3").Value = "Bags"
$ws.Range("E4").Value = "Mugs"
$ws.Range("E5").Value = "Clothing"
$ws.Range("N2").ClearContents()
$ws.Range("N3").ClearContents()

# ---------------------------------------------------------------------------
# 3) Rebuild every hyperlink in the exact order the final workbook uses.
#    (Deleting the whole collection first keeps relationship ids compact and
#    matches the target's observed numbering: rId1..rId13, table part last.)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Test@123", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Test@123")
$ws.Hyperlinks.Add($ws.Range("M3"), "mailto:waseyraby@gmail.com")
$ws.Hyperlinks.Add($ws.Range("M2"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("M5"), "mailto:waseyraby@gmail.com")
$ws.Hyperlinks.Add($ws.Range("M4"), "mailto:waseyraby@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Test@123", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Test@123")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Test@123", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Test@123")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Test@123", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Test@123")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:waseyrabby@btinternet.com")

# ---------------------------------------------------------------------------
# 4) Hyperlinks.Add() can overwrite the displayed cell text (when a
#    TextToDisplay was supplied) and always reassigns the cell style, so
#    restore the real values now that every hyperlink object exists.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 718756
$ws.Range("C3").Value = 718756
$ws.Range("C4").Value = 718756
$ws.Range("C5").Value = 718756
$ws.Range("B2").Value = "waseyrabby@btinternet.com"
$ws.Range("B3").Value = "waseyrabby@btinternet.com"
$ws.Range("B4").Value = "waseyrabby@btinternet.com"
$ws.Range("B5").Value = "waseyrabby@btinternet.com"
$ws.Range("B6").Value = "waseyrabby@btinternet.com"

# ---------------------------------------------------------------------------
# 5) Restore the original cell styles that Hyperlinks.Add() clobbered.
# ---------------------------------------------------------------------------
$ws.Range("N10").Copy()
$ws.Range("B2:B6").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C2:C5").PasteSpecial(-4122)
$ws.Range("M2").PasteSpecial(-4122)
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M5").PasteSpecial(-4122)

# Clear the scratch cells back to their pristine, empty state.
$ws.Range("C10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").ClearContents()

# ---------------------------------------------------------------------------
# 6) Column widths (Excel auto-resized several columns once the longer
#    e-mail address / new category strings were entered).
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 36.6640625
$ws.Columns.Item(11).ColumnWidth = 17
$ws.Columns.Item(12).ColumnWidth = 17.83203125
$ws.Columns.Item(13).ColumnWidth = 29.1640625
$ws.Columns.Item(14).ColumnWidth = 26.33203125

# ---------------------------------------------------------------------------
# 7) Sheet view / selection.
# ---------------------------------------------------------------------------
$ws.Range("F2").Select()
$excel.ActiveWindow.ScrollColumn = 2
